# Ui changes and gitignore
# Updates the Quiz Results sheet: corrects a couple of existing rows and
# appends several new quiz-result rows pulled from the DB export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix rows 5-7 (values shifted slightly in the latest DB pull) ---
$ws.Range("A5").Value = 109
$ws.Range("D5").Value = 40
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 12

$ws.Range("A6").Value = 110
$ws.Range("G6").Value = 23

$ws.Range("A7").Value = 112
$ws.Range("B7").Value = "Deepa"
$ws.Range("C7").Value = "I0796921"
$ws.Range("D7").Value = 60
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 12
$ws.Range("H7").Value = "2025-04-27 23:12:00"

# --- Append new rows 8-14 ---
$newRows = @(
    @(113, "Deepa",    "312424",       100, 5, 5, 12,  "2025-04-27 23:13:00"),
    @(114, "Meeryte",  "I012345566",   20,  5, 1, 24,  "2025-04-27 23:13:00"),
    @(115, "Zumba",    "I5674572",     100, 5, 5, 1,   "2025-04-27 23:13:00"),
    @(116, "ddfghd",   "2124235",      100, 5, 5, 2,   "2025-04-27 23:14:00"),
    @(117, "24242352", "2435325634",   40,  5, 2, 224, "2025-04-27 23:14:00"),
    @(118, "asdgs",    "343463",       40,  5, 2, 12,  "2025-04-27 23:21:00"),
    @(119, "Deepa",    "12415135",     20,  5, 1, 60,  "2025-04-27 23:22:00")
)

$r = 8
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# --- Column C needs to widen slightly to fit the new, longer I-Numbers ---
$ws.Columns.Item(3).ColumnWidth = 12.71
